# Update a set of numeric values in Sheet1 to reflect the refreshed
# "RandomForest" imputation results (Update Name of Algo).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = -21.03559999999998
$ws.Range("A6").Value = -22.64050000000001
$ws.Range("A7").Value = -20.36059999999997
$ws.Range("E7").Value = 15.4779
$ws.Range("E12").Value = 17.7097
$ws.Range("E15").Value = 16.3476
$ws.Range("A16").Value = -21.51889999999998
$ws.Range("A20").Value = -19.86
$ws.Range("E20").Value = 15.98059999999999
$ws.Range("E21").Value = 17.046
$ws.Range("E22").Value = 17.039
$ws.Range("E23").Value = 16.21899999999999
$ws.Range("A28").Value = -22.10690000000001
$ws.Range("A29").Value = -21.43569999999998
$ws.Range("E29").Value = 17.50410000000002
$ws.Range("A32").Value = -21.13089999999999
$ws.Range("E34").Value = 17.32180000000001
$ws.Range("A40").Value = -20.33700000000001
$ws.Range("E42").Value = 16.3679
$ws.Range("E43").Value = 17.435
$ws.Range("E44").Value = 16.62639999999999
$ws.Range("E45").Value = 16.4315
$ws.Range("A46").Value = -21.81880000000001
$ws.Range("E46").Value = 17.13399999999999
$ws.Range("E50").Value = 16.2235
$ws.Range("A51").Value = -21.79439999999999
$ws.Range("E51").Value = 17.21470000000001
$ws.Range("A52").Value = -22.21639999999999
$ws.Range("A57").Value = -22.58780000000002
$ws.Range("A59").Value = -21.9989
$ws.Range("A62").Value = -21.9634
$ws.Range("A66").Value = -21.55060000000001
$ws.Range("E66").Value = 17.02440000000002
$ws.Range("E67").Value = 17.23310000000002
$ws.Range("A73").Value = -20.3498
$ws.Range("A74").Value = -21.68679999999998
$ws.Range("E79").Value = 18.20500000000002
$ws.Range("E84").Value = 16.59369999999999
$ws.Range("A92").Value = -21.34830000000001
$ws.Range("E92").Value = 18.59230000000001
$ws.Range("E97").Value = 16.65110000000001
$ws.Range("A100").Value = -22.06199999999999
